$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.580.29"
$ws.Range("E2").Value = "  -4.85%  "

$ws.Range("D3").Value = "3.256.76"
$ws.Range("E3").Value = "  -8.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.31"
$ws.Range("E5").Value = "  -4.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.27"
$ws.Range("E6").Value = "  -11.84%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.246.50"
$ws.Range("E8").Value = "  -8.28%  "

$ws.Range("E10").Value = "  -12.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.89"
$ws.Range("E11").Value = "  -4.83%  "

$ws.Range("E12").Value = "  -13.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.91"
$ws.Range("E13").Value = "  -16.90%  "

$ws.Range("E14").Value = "  -11.14%  "

$ws.Range("D15").Value = "3.782.46"
$ws.Range("E15").Value = "  -8.12%  "

$ws.Range("D16").Value = "67.668.03"
$ws.Range("E16").Value = "  -4.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "547.47"
$ws.Range("E17").Value = "  -10.72%  "

$ws.Range("D18").Value = "3.258.66"
$ws.Range("E18").Value = "  -8.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.30"
$ws.Range("E19").Value = "  -13.94%  "

$ws.Range("E20").Value = "  -5.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.36"
$ws.Range("E21").Value = "  -13.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.771"
$ws.Range("E22").Value = "  -13.42%  "

$ws.Range("E23").Value = "  -13.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.17"
$ws.Range("E24").Value = "  -12.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.64"
$ws.Range("E25").Value = "  -13.54%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.21"
$ws.Range("E27").Value = "  -15.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.29"
$ws.Range("E28").Value = "  -9.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "29.73"
$ws.Range("E29").Value = "  -12.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  -17.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").Value = "  -9.03%  "

$ws.Range("E32").Value = "  -10.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "557.71"
$ws.Range("E33").Value = "  -11.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.64"
$ws.Range("E34").Value = "  -18.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.82"
$ws.Range("E35").Value = "  -15.67%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.93"
$ws.Range("E37").Value = "  -5.48%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0444"
$ws.Range("E38").Value = "  -7.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.29"
$ws.Range("E39").Value = "  -14.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0855"
$ws.Range("E40").Value = "  -14.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("E41").Value = "  -11.99%  "

$ws.Range("D42").Value = "2.949.45"
$ws.Range("E42").Value = "  -12.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  -24.52%  "

$ws.Range("E44").Value = "  -16.10%  "

$ws.Range("D45").Value = "0.0₃0590"
$ws.Range("E45").Value = "  -20.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.39"
$ws.Range("E46").Value = "  -20.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.46"
$ws.Range("E47").Value = "  -17.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.15"
$ws.Range("E48").Value = "  -16.56%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.50"
$ws.Range("E50").Value = "  -5.87%  "

$ws.Range("E51").Value = "  -12.68%  "
